$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# =========================================================================
# Sheet "NutritionalData" (sheet3): add poki-bowl ingredient nutrition rows
# 160-172, plus totals row 173.
# =========================================================================
$ws3.Range("A160").Value2 = 'salmon, serving 2 pcs sushi,  https://www.calorieking.com/us/en/foods/f/calories-in-japanese-raw-salmon-sushi-nigiri/Q3eMQgVHRfOct_GXeY9aog'
$ws3.Range("B160").Value2 = 134
$ws3.Range("C160").Value2 = 3.9
$ws3.Range("D160").Value2 = 0.9
$ws3.Range("E160").Value2 = 7.1
$ws3.Range("F160").Value2 = 16.7
$ws3.Range("G160").Value2 = 0.6
$ws3.Range("H160").Value2 = 193

$ws3.Range("A161").Value2 = 'ahi tunu, serving 2 pcs sushi, https://www.calorieking.com/us/en/foods/f/calories-in-japanese-raw-tuna-sushi-nigiri/dRH0pXbUTO-_tklypq6u1Q'
$ws3.Range("B161").Value2 = 106
$ws3.Range("C161").Value2 = 0.4
$ws3.Range("D161").Value2 = 0.1
$ws3.Range("E161").Value2 = 8
$ws3.Range("F161").Value2 = 16.7
$ws3.Range("G161").Value2 = 0.6
$ws3.Range("H161").Value2 = 186

$ws3.Range("A162").Value2 = 'spicy mayo, serving 2 oz 1 pkg condiment, https://www.calorieking.com/us/en/foods/f/calories-in-condiments-spicy-mayo/XKs-hi5UQcmDmdyUttVVwg'
$ws3.Range("B162").Value2 = 40
$ws3.Range("C162").Value2 = 4.5
$ws3.Range("D162").Value2 = 1
$ws3.Range("E162").Value2 = 0
$ws3.Range("F162").Value2 = 0
$ws3.Range("G162").Value2 = 0
$ws3.Range("H162").Value2 = 40

$ws3.Range("A163").Value2 = 'teriyaki sauce, serving 1 tbs 0.6 oz, https://www.calorieking.com/us/en/foods/f/calories-in-sauces-teriyaki-sauce/31p8AsaFT7q7WEs1X9AO1A'
$ws3.Range("B163").Value2 = 16
$ws3.Range("C163").Value2 = 0.1
$ws3.Range("D163").Value2 = 0
$ws3.Range("E163").Value2 = 1.1
$ws3.Range("F163").Value2 = 2.8
$ws3.Range("G163").Value2 = 0.1
$ws3.Range("H163").Value2 = 690

$ws3.Range("A164").Value2 = 'pineapple_1cup, 5.8 oz, https://www.calorieking.com/us/en/foods/f/calories-in-fresh-fruits-pineapple-raw/m4JttugzRT2KyGMkrrZkzQ'
$ws3.Range("B164").Value2 = 82
$ws3.Range("C164").Value2 = 0.2
$ws3.Range("D164").Value2 = 0.1
$ws3.Range("E164").Value2 = 0.9
$ws3.Range("F164").Value2 = 21.6
$ws3.Range("G164").Value2 = 2.3
$ws3.Range("H164").Value2 = 2

$ws3.Range("A165").Value2 = 'pinapple poki bowl 1/4 cup'
$ws3.Range("B165").Formula = '=B164/4'
$ws3.Range("C165").Formula = '=C164/4'
$ws3.Range("D165").Formula = '=D164/4'
$ws3.Range("E165").Formula = '=E164/4'
$ws3.Range("F165").Formula = '=F164/4'
$ws3.Range("G165").Formula = '=G164/4'
$ws3.Range("H165").Formula = '=H164/4'

$ws3.Range("A166").Value2 = 'brown rice, 1 cup 6.9 oz, https://www.calorieking.com/us/en/foods/f/calories-in-rice-long-grain-brown-rice-cooked/_W12T-lpSI-Nm25pKDkvrw'
$ws3.Range("B166").Value2 = 216
$ws3.Range("C166").Value2 = 1.8
$ws3.Range("D166").Value2 = 0.4
$ws3.Range("E166").Value2 = 5
$ws3.Range("F166").Value2 = 44.8
$ws3.Range("G166").Value2 = 3.5
$ws3.Range("H166").Value2 = 10

$ws3.Range("A167").Value2 = 'poki brown rice 1 1/2 cups'
$ws3.Range("B167").Formula = '=B166*3/2'
$ws3.Range("C167").Formula = '=C166*3/2'
$ws3.Range("D167").Formula = '=D166*3/2'
$ws3.Range("E167").Formula = '=E166*3/2'
$ws3.Range("F167").Formula = '=F166*3/2'
$ws3.Range("G167").Formula = '=G166*3/2'
$ws3.Range("H167").Formula = '=H166*3/2'

$ws3.Range("A168").Value2 = 'cucumbers, serving 1 cup 3.7 oz, https://www.calorieking.com/us/en/foods/f/calories-in-fresh-or-dried-vegetables-cucumber-with-peel-raw/nKvHV3A9TamX5TW2U8mjdA'
$ws3.Range("B168").Value2 = 16
$ws3.Range("C168").Value2 = 0.1
$ws3.Range("D168").Value2 = 0.1
$ws3.Range("E168").Value2 = 0.7
$ws3.Range("F168").Value2 = 3.8
$ws3.Range("G168").Value2 = 0.5
$ws3.Range("H168").Value2 = 2

$ws3.Range("A169").Value2 = 'cucumbers poki 1/4 cup '
$ws3.Range("B169").Formula = '=B168/4'
$ws3.Range("C169").Formula = '=C168/4'
$ws3.Range("D169").Formula = '=D168/4'
$ws3.Range("E169").Formula = '=E168/4'
$ws3.Range("F169").Formula = '=F168/4'
$ws3.Range("G169").Formula = '=G168/4'
$ws3.Range("H169").Formula = '=H168/4'

$ws3.Range("A170").Value2 = 'sesame seeds serving 1 tbs 0.3 oz, https://www.calorieking.com/us/en/foods/f/calories-in-seeds-whole-sesame-seeds-roasted-toasted/FyvsC9v-QiOYvcWSJCvszw'
$ws3.Range("B170").Value2 = 51
$ws3.Range("C170").Value2 = 4.3
$ws3.Range("D170").Value2 = 0.6
$ws3.Range("E170").Value2 = 1.5
$ws3.Range("F170").Value2 = 2.3
$ws3.Range("G170").Value2 = 1.3
$ws3.Range("H170").Value2 = 1

$ws3.Range("A171").Value2 = 'ginger root poki bowl 1 teaspoon 0.1 oz seving, https://www.calorieking.com/us/en/foods/f/calories-in-fresh-or-dried-vegetables-ginger-root-raw/-TXI3lP0RCav-6oUYlV6Ww'
$ws3.Range("B171").Value2 = 2
$ws3.Range("C171").Value2 = 0.1
$ws3.Range("D171").Value2 = 0.1
$ws3.Range("E171").Value2 = 0.1
$ws3.Range("F171").Value2 = 0.4
$ws3.Range("G171").Value2 = 0.1
$ws3.Range("H171").Value2 = 1

$ws3.Range("A172").Value2 = 'wasabi 1 teaspoon serving, https://www.calorieking.com/us/en/foods/f/calories-in-japanese-wasabi-root-raw/hcWfFF5TRE63dEPlZbo-iA'
$ws3.Range("B172").Value2 = 5
$ws3.Range("C172").Value2 = 0.1
$ws3.Range("D172").Value2 = 0
$ws3.Range("E172").Value2 = 0.2
$ws3.Range("F172").Value2 = 1.2
$ws3.Range("G172").Value2 = 0.4
$ws3.Range("H172").Value2 = 1

# Row 173: totals (2x the 1-tbs spicy mayo baseline row 162)
$ws3.Range("B173").Formula = '=B162*2'
$ws3.Range("C173").Formula = '=C162*2'
$ws3.Range("D173").Formula = '=D162*2'
$ws3.Range("E173").Formula = '=E162*2'
$ws3.Range("F173").Formula = '=F162*2'
$ws3.Range("G173").Formula = '=G162*2'
$ws3.Range("H173").Formula = '=H162*2'


# =========================================================================
# Sheet "researchMeasures" (sheet1): row 66 gets the poki-bowl lunch entry
# written up (note column + nutrition totals), plus a couple of small
# value tweaks; row 67 gains a formatted-but-empty AR cell.
# =========================================================================
$ws1.Range("I66").Value2 = 0

$ws1.Range("Z66").Value2 = 'Woke up at 530 am by alarm, reviewed the discussion and power point slides in genetics and folded laundry from yesterday after cleaning up one pet mess, making my coffee, and feeding the babies, had 2nd cup of coffee after folding the laundry. I bought some French press from IKEA when I went earlier in the week and had some coffee grounds from Starbucks ground French press style, but haven''t had time to clean the item and follow instructions for cold brew of french press coffee. Had a lg BM after 2nd cup of coffee. Then took measurements crunched for time. Had a banana, frozen blueberry, pea protein and cocoa powder smoothie for breakfast with a little bit of coffee for 3rd cup in it. For lunch I had a poki bowl with 1 scoop of salmon, one of the tuna, the spicy crab meat that I didn''t really eat, over brown rice about 1 1/2 cups with spicy mayo 2 tbs and teriyaki sauce about 2 tbs with ginger 1 tbs, pineapples 1 scoop about 1 slice chopped, sesame seeds toasted about 1/2 tbs, about 1/4 cup chopped cucumbers too. And my 4th cup of coffee from work''s keurig the Donut Shop one today. After work ate 3 quesadillas normal style and shared with the babies, later had a large piece of chocolate cake about 2 cupcakes worth in size. My amazon order arrived of the duffel bag after doing some writing of nomenclature for each individual compound from the week 3 worksheet with polyatomic charges and fixed or variable metal charges as well as other notes on notecards. Then completed the week 5 power point fill in for genetics. We have a chapter 4 quiz in chemistry Tuesday, homework due Wednesday but the masteringLab was froze earlier and couldn''t look at the 3 hours they estimate to complete it, also a lab dry lab worksheet on ionization and oxidation from the chapter 5 slides I need to go over and make flashcards for. I need to get more flashcards and maybe some of those little notecard boxes they sell at Michael''s craft store by my work. My ankle is still swollen and my right knee. I wanted to workout but I didn''t due to catching up on prepping for stuying and studying while preparing the note cards. I got some plastic bags for laundered linens after watching them in gallon size and the sandwhich bag size for food or small work items for private clients and dryer sheets when picking up one of the other robes I ordered at the hub spot. And another one is waiting for me to pick it up there tomorrow morning before work, that and another little basket. Not bad size. Probably sell the same thing at the dollar tree but convenient to order it on Amazon for $10 and use my monthly prime subscription for free delivery. Still need the biofreeze order and especially the hydrocollator that hasn''t shipped yet. I don''t plan on using it any time soon, but want it here and ready. Tired. Been tired. Have a discussion I need to do before the end of the weekend. Rag my 3rd day and only medium. yesterday was medium light, and the first day was light. Usually its spotty-light, med-light, heavy, med-heavy, med-light, light for 6 days approximately. It might be the red fennel pasta and not eating it with cheese and having the peppers because of the phytohormones in fennel. I had my 5th cup of coffee while eating the quesadillas and making the flashcards for chemistry on nomenclature.'

$ws1.Range("AA66").Value2 = "1 serving pea protein`n(120`t2`t0`t18`t6`t1`t360)`n2 tbs cocoa`n(20`t1`t0`t2`t6`t2`t0)`n1/2 cup frozen blueberries`n(42`t0`t0`t1`t13`t2`t1)`n1 banana`n(105`t0`t0`t1`t27`t3`t1)`nPoki Bowl`n1/4 cup salmon `n(134`t3.9`t0.9`t7.1`t16.7`t0.6`t193)`n1/4 cup tuna bluefin ahi`n(106`t0.4`t0.1`t8`t16.7`t0.6`t186)`n2 tbs spicy mayo`n(80`t9`t2`t0`t0`t0`t80)`n2 tbs teriyaki sauce`n(32`t0.2`t0`t2.2`t5.6`t0.2`t1380)`n1/4 cup pineapple`n(20.5`t0.05`t0.025`t0.225`t5.4`t0.575`t0.5)`n1/4 cup cucumbers`n(4`t0.025`t0.025`t0.175`t0.95`t0.125`t0.5)`n1 tbs sesame seeds`n(51`t4.3`t0.6`t1.5`t2.3`t1.3`t1)`n1 1/2 cups brown rice`n(324`t2.7`t0.6`t7.5`t67.2`t5.25`t15)`n2 tbs ginger`n(2`t0.1`t0.1`t0.1`t0.4`t0.1`t1)`n1/4 tbs wasabi`n(5`t0.1`t0`t0.2`t1.2`t0.4`t1)`n6 corn tortillas Guerrero brand`n(300`t3`t0`t6`t63`t6`t60)`n1/2 cup mozzarella winco brand`n(160`t10`t7`t12`t2`t0`t380)`n2 chocolate cupcakes`n(164`t5.4`t1.2`t1.7`t29.2`t0.7`t176)`n=120+20+42+105+134+106+80+32+20.5+4+51+324+2+5+300+160+164`n=2+1+0+0+3.9+0.4+9+0.2+0.005+0.025+4.3+2.7+0.1+0.1+3+10+5.4`n=0+0+0+0+0.9+0.1+2+0+0.025+0.025+0.6+0.6+0.1+0+0+7+1.2`n=18+2+1+1+7.1+8+0+2.2+0.225+0.175+1.5+7.5+0.1+0.2+6+12+1.7`n=6+6+13+27+16.7+16.7+0+5.6+5.4+0.95+2.3+67.2+0.4+1.2+63+2+29.2`n=1+2+2+3+0.6+0.6+0.2+0.575+0.125+1.3+5.25+0.1+0.4+6+0+0.7`n=360+0+1+1+193+186+80+1380+0.5+0.5+1+15+1+1+60+380+176`n"
$ws1.Range("AA66").WrapText = $true

$ws1.Range("AB66").Formula = '=120+20+42+105+134+106+80+32+20.5+4+51+324+2+5+300+160+164'
$ws1.Range("AC66").Formula = '=2+1+0+0+3.9+0.4+9+0.2+0.005+0.025+4.3+2.7+0.1+0.1+3+10+5.4'
$ws1.Range("AD66").Formula = '=0+0+0+0+0.9+0.1+2+0+0.025+0.025+0.6+0.6+0.1+0+0+7+1.2'
$ws1.Range("AE66").Formula = '=18+2+1+1+7.1+8+0+2.2+0.225+0.175+1.5+7.5+0.1+0.2+6+12+1.7'
$ws1.Range("AF66").Formula = '=6+6+13+27+16.7+16.7+0+5.6+5.4+0.95+2.3+67.2+0.4+1.2+63+2+29.2'
$ws1.Range("AG66").Formula = '=1+2+2+3+0.6+0.6+0.2+0.575+0.125+1.3+5.25+0.1+0.4+6+0+0.7'
$ws1.Range("AH66").Formula = '=360+0+1+1+193+186+80+1380+0.5+0.5+1+15+1+1+60+380+176'

# Extend the AI:AN ratio (nutrient / calories) shared formulas down to row 66
$ws1.Range('AI66').Formula = '=$AC66/$AB66'
$ws1.Range('AJ66').Formula = '=$AD66/$AB66'
$ws1.Range('AK66').Formula = '=$AE66/$AB66'
$ws1.Range('AL66').Formula = '=$AF66/$AB66'
$ws1.Range('AM66').Formula = '=$AG66/$AB66'
$ws1.Range('AN66').Formula = '=$AH66/$AB66'

$ws1.Range("AO66").Value2 = 5
$ws1.Range("AR66").Value2 = 0

# Row 67 gains a blank, wrap-text-formatted AR cell (matches column AR style)
$ws1.Range("AR67").WrapText = $true

# =========================================================================
# Window / view state clean-up so the workbook re-opens scrolled to the
# areas that were being edited (cosmetic, matches authors final screen).
# =========================================================================
$ws3.Activate()
$excel.Goto($ws3.Range("A52"), $false)
$ws3.Range("B56:H56").Select()

$ws1.Activate()
$excel.Goto($ws1.Range("A60"), $false)
$ws1.Range("O1").Select()
$ws1.Range("AO68").Select()
